$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on both the "展览" sheet and
# the "全部类型" sheet, which mirror each other's data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 240
    $ws.Range("F3").Value = 361
}
